$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Big boy"
$ws.Range("C9").Value = "Biggest boy arouund"
$ws.Range("D9").Value = "Canada"
$ws.Range("E9").Value = "2021-10-11T13:31:41.740Z"
$ws.Range("F9").Value = "Sent"
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = "AD"
$ws.Range("I9").Value = "los santos"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "My lovely Job"
$ws.Range("C10").Value = "Tiny miney einie weenie time town"
$ws.Range("D10").Value = "Google"
$ws.Range("E10").Value = "2021-10-11T13:40:00.771Z"
$ws.Range("F10").Value = "Sent"
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = "AI"
$ws.Range("I10").Value = "asdas"

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "s"
$ws.Range("C11").Value = "d  dasdasd"
$ws.Range("D11").Value = "asdad"
$ws.Range("E11").Value = "2021-10-11T13:44:03.679Z"
$ws.Range("F11").Value = "Sent"
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = "Angola"
$ws.Range("I11").Value = "ASD"
